# BL Audit Form - "update report" sheet - 11-11-24 report update
# Applies the numeric/text changes captured in the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("update report")
$ws.Activate()

# --- Report date (B1) and payment note (F34) -----------------------------
# Both hold date-ish text that Excel would otherwise auto-convert to a date
# serial if assigned directly. Write them as a text-producing formula first,
# then Copy / PasteSpecial(values) them onto themselves - this keeps the
# original cell style (no NumberFormat change) while storing a literal
# shared string instead of a date number.
# F34 is processed first so the shared-string table ends up ordered the
# same way as in the target file (F34's string lands before B1's).

$ws.Range("F34").Formula = "=""11.11.2024 payment """
$ws.Range("F34").Copy()
$ws.Range("F34").PasteSpecial(-4163)

$ws.Range("B1").Formula = "=""11.11.2024"""
$ws.Range("B1").Copy()
$ws.Range("B1").PasteSpecial(-4163)

$excel.CutCopyMode = 0

# --- Stock section (rows 9-20) quantities ---------------------------------
$ws.Range("C9").Value = 311307
$ws.Range("C10").Value = 60
$ws.Range("C11").Value = 150
$ws.Range("C12").Value = 7120
$ws.Range("C16").Value = 17
$ws.Range("C17").Value = 76
$ws.Range("C19").Value = 80
$ws.Range("C20").ClearContents()

# --- Dues / liabilities ----------------------------------------------------
$ws.Range("E22").Value = 44292
$ws.Range("E23").Value = 83702

# --- Bank guarantee / commission section -----------------------------------
$ws.Range("E27").Value = 37750
$ws.Range("E29").Value = 12474

# --- Credit section ----------------------------------------------------
$ws.Range("E34").Value = 100000

# --- View state: selected cell moved to E35, scrolled further down --------
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E35").Select()
